$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users Stories")

# --- New shared-string texts (must be created in this exact order so the
# underlying shared string table gets indices 51..56 in the same order the
# author's Excel produced them) -----------------------------------------
$s51 = "Como Personal`nQuero listar as Avaliações Físicas de um aluno`nPara que possa visualizalas"
$s52 = "Como Personal`nQuero listar as Sereies de um aluno`nPara que possa visualizalas"
$s53 = "Como Personal`nQuero listar os Alunos`nPara que possa visualizalos"
$s54 = "Como Personal`nQuero listar os Exercícios`nPara que possa Selecionalos"
$s55 = "Como Aluno`nQuero listar as Minhas Séries`nPara que possa Visualizalas"
$s56 = "Como Aluno`nQuero listar as minhas Avaliações Físicas `nPara que possa visualizalas"

# --- Row 35: new section header "User Story Funcional" (same look as the
# header rows above it, e.g. row 33 "User Story do Sistema"). The range is
# merged *before* the formatting is copied in, otherwise Excel recomputes
# the merged region's border and produces new (duplicate) border/style
# records instead of reusing the existing header style (66/67/68). -------
$ws.Range("B35:H35").Merge()
$ws.Range("B33:H33").Copy()
$ws.Range("B35:H35").PasteSpecial(-4122)
$ws.Range("B35").Value = "User Story Funcional"

# --- Rows 36-41: new user stories, formatted like row 34 (the previous
# last data row in the "User Story Funcional" table) ---------------------
$ws.Range("B34:H34").Copy($ws.Range("B36:H36"))
$ws.Range("B34:H34").Copy($ws.Range("B37:H37"))
$ws.Range("B34:H34").Copy($ws.Range("B38:H38"))
$ws.Range("B34:H34").Copy($ws.Range("B39:H39"))
$ws.Range("B34:H34").Copy($ws.Range("B40:H40"))
$ws.Range("B34:H34").Copy($ws.Range("B41:H41"))

# NOTE: new shared strings get appended to the workbook's shared-string
# table in the order the *cell values are actually assigned* (not the
# order the PowerShell variables were declared). To reproduce the exact
# target table order (51=s51, 52=s52, ... 56=s56) the cells must receive
# their text in that same sequence: C39, C38, C36, C37, C40, C41.
$ws.Range("C39").Value = $s51
$ws.Range("C38").Value = $s52
$ws.Range("C36").Value = $s53
$ws.Range("C37").Value = $s54
$ws.Range("C40").Value = $s55
$ws.Range("C41").Value = $s56

$ws.Range("B36").Value = 29
$ws.Range("D36").Value = 4
$ws.Range("E36").Value = 3
$ws.Range("F36").Value = 3
$ws.Range("G36").Value = 3
$ws.Range("H36").Value = ""

$ws.Range("B37").Value = 30
$ws.Range("D37").Value = 4
$ws.Range("E37").Value = 3
$ws.Range("F37").Value = 3
$ws.Range("G37").Value = 3
$ws.Range("H37").Value = ""

$ws.Range("B38").Value = 31
$ws.Range("D38").Value = 4
$ws.Range("E38").Value = 3
$ws.Range("F38").Value = 3
$ws.Range("G38").Value = 3
$ws.Range("H38").Value = ""

$ws.Range("B39").Value = 32
$ws.Range("D39").Value = 4
$ws.Range("E39").Value = 3
$ws.Range("F39").Value = 3
$ws.Range("G39").Value = 3
$ws.Range("H39").Value = ""

$ws.Range("B40").Value = 33
$ws.Range("D40").Value = 4
$ws.Range("E40").Value = 3
$ws.Range("F40").Value = 3
$ws.Range("G40").Value = 3
$ws.Range("H40").Value = ""

$ws.Range("B41").Value = 34
$ws.Range("D41").Value = 4
$ws.Range("E41").Value = 3
$ws.Range("F41").Value = 3
$ws.Range("G41").Value = 3
$ws.Range("H41").Value = ""
